$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $c = $ws.Range($cell)
    $c.Value = "'" + $value
    $c.Style = "Normal"
}

Set-TextCell "D2" "27.122.68"
Set-TextCell "E2" "  -0.31%  "
Set-TextCell "D3" "1.892.82"
Set-TextCell "E3" "  -0.73%  "
Set-TextCell "E4" "  +0.15%  "
Set-TextCell "D5" "307.14"
Set-TextCell "E5" "  -0.29%  "
Set-TextCell "E6" "  +0.10%  "
Set-TextCell "D7" "0.5228"
Set-TextCell "E7" "  -0.49%  "
Set-TextCell "D8" "0.3761"
Set-TextCell "E8" "  -0.63%  "
Set-TextCell "D9" "0.07265"
Set-TextCell "E9" "  -0.19%  "
Set-TextCell "D10" "21.09"
Set-TextCell "E10" "  -0.77%  "
Set-TextCell "D11" "0.9003"
Set-TextCell "E11" "  +0.27%  "
Set-TextCell "D12" "0.08179"
Set-TextCell "E12" "  +6.37%  "
Set-TextCell "D13" "1.935.32"
Set-TextCell "D14" "96.15"
Set-TextCell "E14" "  +1.14%  "
Set-TextCell "D15" "5.289"
Set-TextCell "E15" "  +0.22%  "
Set-TextCell "D16" "1.002"
Set-TextCell "E16" "  +0.14%  "
Set-TextCell "D17" "0.000008574"
Set-TextCell "E17" "  -0.60%  "
Set-TextCell "E18" "  +0.55%  "
Set-TextCell "E19" "  +0.15%  "
Set-TextCell "D20" "27.141.89"
Set-TextCell "E20" "  -0.46%  "
Set-TextCell "D21" "5.083"
Set-TextCell "E21" "  +0.19%  "
Set-TextCell "D22" "10.69"
Set-TextCell "E22" "  +0.49%  "
Set-TextCell "D23" "6.406"
Set-TextCell "E23" "  -0.55%  "
Set-TextCell "D24" "148.06"
Set-TextCell "E24" "  +1.54%  "
Set-TextCell "D25" "2.287"
Set-TextCell "E25" "  -1.26%  "
Set-TextCell "B26" "Toncoin"
Set-TextCell "C26" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D26" "1.741"
Set-TextCell "E26" "  +0.31%  "
Set-TextCell "B27" "EthereumClassic"
Set-TextCell "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D27" "18.16"
Set-TextCell "E27" "  +0.06%  "
Set-TextCell "D28" "115.02"
Set-TextCell "E28" "  +0.17%  "
Set-TextCell "D29" "4.790"
Set-TextCell "E29" "  -0.59%  "
Set-TextCell "D30" "4.836"
Set-TextCell "E30" "  -2.64%  "
Set-TextCell "D31" "0.09225"
Set-TextCell "E31" "  +0.02%  "
Set-TextCell "D32" "0.05034"
Set-TextCell "E32" "  -0.66%  "
Set-TextCell "D33" "0.7877"
Set-TextCell "E33" "  -3.60%  "
Set-TextCell "D34" "1.214"
Set-TextCell "E34" "  -2.24%  "
Set-TextCell "D35" "3.422"
Set-TextCell "E35" "  +3.36%  "
Set-TextCell "D36" "2.958"
Set-TextCell "E36" "  -1.27%  "
Set-TextCell "D37" "2.594"
Set-TextCell "E37" "  -0.03%  "
Set-TextCell "D38" "0.5698"
Set-TextCell "E38" "  +0.18%  "
Set-TextCell "D39" "0.01983"
Set-TextCell "E39" "  -0.43%  "
Set-TextCell "E40" "  -0.08%  "
Set-TextCell "E41" "  +0.25%  "
Set-TextCell "E42" "  -1.23%  "
Set-TextCell "D43" "116.42"
Set-TextCell "E43" "  -2.41%  "
Set-TextCell "D44" "0.1516"
Set-TextCell "E44" "  -0.06%  "
Set-TextCell "D45" "0.4854"
Set-TextCell "E45" "  +0.30%  "
Set-TextCell "E46" "  +0.13%  "
Set-TextCell "D47" "10.03"
Set-TextCell "E47" "  -1.64%  "
Set-TextCell "D48" "1.620"
Set-TextCell "E48" "  -0.35%  "
Set-TextCell "D49" "38.10"
Set-TextCell "E49" "  +1.33%  "
Set-TextCell "D50" "63.49"
Set-TextCell "E50" "  -0.37%  "
Set-TextCell "D51" "0.05935"
Set-TextCell "E51" "  +0.03%  "
